$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.702.61'
$ws.Range('E2').Value = '  +0.87%  '

$ws.Range('D3').Value = '1.816.07'
$ws.Range('E3').Value = '  +1.12%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.83%  '

$ws.Range('E6').Value = '  +2.34%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '35.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.301'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0698'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.69%  '

$ws.Range('E11').Value = '  +0.38%  '

$ws.Range('D12').Value = '2.076.34'
$ws.Range('E12').Value = '  +1.03%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.09%  '

$ws.Range('D14').Value = '1.824.12'
$ws.Range('E14').Value = '  +1.65%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.646'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.55%  '

$ws.Range('D16').Value = '34.717.23'
$ws.Range('E16').Value = '  +1.01%  '

$ws.Range('E17').Value = '  +3.15%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.99%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.71%  '

$ws.Range('D20').Value = '0.0₃0804'
$ws.Range('E20').Value = '  +0.05%  '

$ws.Range('E21').Value = '  +5.02%  '

$ws.Range('E22').Value = '  +0.06%  '

$ws.Range('E23').Value = '  +0.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '171.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.27%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.75%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.79'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.37%  '

$ws.Range('E28').Value = '  +1.45%  '

$ws.Range('E29').Value = '  -0.21%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.58%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0533'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.23%  '

$ws.Range('E32').Value = '  +2.12%  '

$ws.Range('E33').Value = '  +1.54%  '

$ws.Range('E34').Value = '  +1.75%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.64'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.93%  '

$ws.Range('D36').Value = '1.420.99'
$ws.Range('E36').Value = '  -1.44%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.686'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('E38').Value = '  +1.43%  '

$ws.Range('E39').Value = '  +0.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.19%  '

$ws.Range('E43').Value = '  +0.42%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.03%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.11'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.37%  '

$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0521'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.59%  '

$ws.Range('E47').Value = '  +0.44%  '

$ws.Range('D48').Value = '1.977.70'
$ws.Range('E48').Value = '  +1.36%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.15%  '

$ws.Range('E50').Value = '  +0.50%  '

$ws.Range('E51').Value = '  -0.01%  '
